$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rebuild the "Products" sheet: new "Category" and "Specs (JSON)" columns,
#    refreshed sample data (2 rows instead of 10), and wider columns.
# ---------------------------------------------------------------------------
$products = $wb.Worksheets.Item(1)

# Wipe the old 7-column / 11-row sample data before laying out the new shape.
$products.Cells.Clear()

# Column widths (character units) -> matches the authored template widths.
$products.Columns.Item(1).ColumnWidth = 40
$products.Columns.Item(2).ColumnWidth = 18
$products.Columns.Item(3).ColumnWidth = 14
$products.Columns.Item(4).ColumnWidth = 20
$products.Columns.Item(5).ColumnWidth = 18
$products.Columns.Item(6).ColumnWidth = 14
$products.Columns.Item(7).ColumnWidth = 14
$products.Columns.Item(8).ColumnWidth = 60
$products.Columns.Item(9).ColumnWidth = 60

# Header row.
$products.Range("A1").Value = "Product Name"
$products.Range("B1").Value = "Category"
$products.Range("C1").Value = "Brand"
$products.Range("D1").Value = "SKU"
$products.Range("E1").Value = "Base Price (MWK)"
$products.Range("F1").Value = "Stock Quantity"
$products.Range("G1").Value = "Condition"
$products.Range("H1").Value = "Description"
$products.Range("I1").Value = "Specs (JSON)"

# Sample row 2.
$products.Range("A2").Value = "iPhone 15 Pro Max 256GB"
$products.Range("B2").Value = "Smartphones"
$products.Range("C2").Value = "Apple"
$products.Range("D2").Value = "IP15PM-256-BLK"
$products.Range("E2").Value = 1500000
$products.Range("F2").Value = 10
$products.Range("G2").Value = "NEW"
$products.Range("H2").Value = "Brand new, sealed in box. 1 year warranty."
$products.Range("I2").Value = '{"storage":"256GB","color":"Black Titanium","ram":"8GB"}'

# Sample row 3.
$products.Range("A3").Value = "Samsung Galaxy S24 Ultra"
$products.Range("B3").Value = "Smartphones"
$products.Range("C3").Value = "Samsung"
$products.Range("D3").Value = "SGS24U-512-GRY"
$products.Range("E3").Value = 1350000
$products.Range("F3").Value = 5
$products.Range("G3").Value = "NEW"
$products.Range("H3").Value = "Factory unlocked. Includes S Pen."
$products.Range("I3").Value = '{"storage":"512GB","color":"Titanium Gray","ram":"12GB"}'

# Keep "numbers stored as text" warnings suppressed over the new extent.
$products.Range("A1:I3").Errors.Item(1).Ignore = $true

# ---------------------------------------------------------------------------
# 2. Add a new "Instructions" sheet right after "Products" with usage notes.
# ---------------------------------------------------------------------------
$instructions = $wb.Worksheets.Add($null, $products)
$instructions.Name = "Instructions"

$instructions.Columns.Item(1).ColumnWidth = 120

$instructions.Range("A1").Value = "BULK UPLOAD INSTRUCTIONS"
$instructions.Range("A2").Value = ""
$instructions.Range("A3").Value = "Required Columns:"
$instructions.Range("A4").Value = "- Product Name: The name of the product"
$instructions.Range("A5").Value = "- Base Price (MWK): Your selling price BEFORE platform fees"
$instructions.Range("A6").Value = "- Stock Quantity: Number of items in stock"
$instructions.Range("A7").Value = ""
$instructions.Range("A8").Value = "Optional Columns:"
$instructions.Range("A9").Value = "- Category, Brand, SKU, Condition, Description, Specs (JSON)"
$instructions.Range("A10").Value = ""
$instructions.Range("A11").Value = "Notes:"
$instructions.Range("A12").Value = "- Remove sample rows before uploading real products"
$instructions.Range("A13").Value = "- Max 200 products per upload"
$instructions.Range("A14").Value = "- Prices are in MWK"

$instructions.Range("A1:A14").Errors.Item(1).Ignore = $true

# Leave the user back on the primary "Products" sheet.
$products.Activate()
